$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D and E, shifting existing D:K data to F:M
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formats from the (now-shifted) old D:E columns (now F:G) into
# the two newly inserted blank columns so they match the rest of the data range
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New quarterly data for columns D (2018-12-31) and E (2018-09-30)
$data = @(
    @{R=7; D=43465; E=43373},
    @{R=8; D="NA"; E=41600},
    @{R=9; D=12200; E=11600},
    @{R=10; D="NA"; E=30000},
    @{R=12; D="NA"; E="NA"},
    @{R=13; D=0; E=0},
    @{R=14; D=0; E=0},
    @{R=15; D=14500; E=13800},
    @{R=17; D="NA"; E=51800},
    @{R=18; D="NA"; E=-10200},
    @{R=20; D="NA"; E=6100},
    @{R=21; D="NA"; E=9700},
    @{R=22; D=20100; E=19300},
    @{R=23; D="NA"; E=-23400},
    @{R=24; D=0; E=0},
    @{R=25; D=0; E=0},
    @{R=26; D=-16400; E=-23400},
    @{R=27; D=-16400; E=-22200},
    @{R=28; D=0; E=0},
    @{R=29; D=0; E=0},
    @{R=30; D=0; E=0},
    @{R=31; D=0; E=0},
    @{R=32; D="NA"; E=-6100},
    @{R=33; D=-16400; E=-22200},
    @{R=34; D=0; E=0},
    @{R=35; D=-16400; E=-22200},
    @{R=38; D=43465; E=43373},
    @{R=41; D=46600; E=22700},
    @{R=42; D="NA"; E=23400},
    @{R=43; D=44000; E=29400},
    @{R=44; D=0; E=0},
    @{R=45; D=13300; E=15000},
    @{R=46; D=104000; E=90500},
    @{R=47; D=69400; E=168900},
    @{R=48; D=1383600; E=1163900},
    @{R=49; D=38200; E=1600},
    @{R=50; D=0; E=0},
    @{R=51; D=0; E=0},
    @{R=52; D=32200; E=31500},
    @{R=53; D=0; E=0},
    @{R=54; D=1627400; E=1456500},
    @{R=57; D=12600; E=8900},
    @{R=58; D="NA"; E=80800},
    @{R=59; D=28600; E=47100},
    @{R=60; D=41200; E=136800},
    @{R=61; D=1205800; E=935500},
    @{R=62; D=0; E=6200},
    @{R=63; D=0; E=0},
    @{R=64; D=0; E=0},
    @{R=65; D=0; E=0},
    @{R=66; D=1247000; E=1078500},
    @{R=68; D=0; E=0},
    @{R=69; D=0; E=0},
    @{R=70; D=0; E=0},
    @{R=71; D=0; E=0},
    @{R=72; D="NA"; E=-125600},
    @{R=73; D=0; E=0},
    @{R=74; D=0; E=0},
    @{R=75; D=0; E=0},
    @{R=76; D=380400; E=378000},
    @{R=77; D=0; E=0},
    @{R=80; D=43465; E=43373},
    @{R=81; D=-16400; E=-22200},
    @{R=83; D=14500; E=13800},
    @{R=84; D=0; E=0},
    @{R=85; D=0; E=0},
    @{R=86; D=0; E=0},
    @{R=87; D=0; E=0},
    @{R=88; D=0; E=0},
    @{R=89; D=-14900; E=7100},
    @{R=91; D=0; E=0},
    @{R=92; D=0; E=0},
    @{R=93; D=0; E=0},
    @{R=94; D=26800; E=1600},
    @{R=96; D=-3000; E=-3000},
    @{R=97; D=0; E=0},
    @{R=98; D=0; E=0},
    @{R=99; D=0; E=0},
    @{R=100; D=-13700; E=-12500},
    @{R=101; D=0; E=0},
    @{R=102; D=-1700; E=-3800}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.R, 4).Value = $item.D
    $ws.Cells.Item($item.R, 5).Value = $item.E
}
